# 6.a.1.1.xlsx — add 2020-2023 data columns (X:AA) and hide the now
# unused 2000-2007 columns (D:K).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New year headers (row 4) ------------------------------------------------
$ws.Range("X4").Value = 2020
$ws.Range("Y4").Value = 2021
$ws.Range("Z4").Value = 2022
$ws.Range("AA4").Value = 2023

# --- New data (rows 5-7) -----------------------------------------------------
# Row 5: "Investment loan"
$ws.Range("X5").Value = 23780
$ws.Range("Y5").Value = 44660
$ws.Range("Z5").Value = 25000
$ws.Range("AA5").Value = 13010

# Row 6: "Investment grant"
$ws.Range("X6").Value = 38240
$ws.Range("Y6").Value = 7950
$ws.Range("Z6").Value = 23000
$ws.Range("AA6").Value = 16390

# Row 7: "Investment loan and grant"
$ws.Range("X7").Value = 62020
$ws.Range("Y7").Value = 52610
$ws.Range("Z7").Value = 48000
$ws.Range("AA7").Value = 29400

# --- Copy formatting from the previous year column (W) onto the new ones ----
$ws.Range("W4:W7").Copy()
$ws.Range("X4:AA7").PasteSpecial(-4122) # xlPasteFormats

# --- Hide the oldest (now empty of new-data relevance) year columns D:K -----
$ws.Range("D1:K1").EntireColumn.Hidden = $true

# --- Row heights: header/band rows grow slightly to fit four more columns ---
$ws.Rows("4").RowHeight = 16.5
$ws.Rows("5").RowHeight = 16.5
$ws.Rows("6").RowHeight = 16.5
$ws.Rows("7").RowHeight = 16.5

$wb.Save()
